$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.899.72'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '1.546.42'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.36'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0582'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('D12').Value = '1.766.83'
$ws.Range('E12').Value = '  -1.15%  '
$ws.Range('D13').Value = '1.546.22'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.512'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Value = '26.891.24'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.40'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '214.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0683'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('E22').Value = '  -2.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('E24').Value = '  -3.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.72'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.103'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.33%  '
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('E31').Value = '  -1.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.91%  '
$ws.Range('D33').Value = '1.356.05'
$ws.Range('E33').Value = '  -3.38%  '
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.962'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.86%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.804'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.06%  '
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.58'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.991'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('D47').Value = '1.681.35'
$ws.Range('E47').Value = '  -1.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.20'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.85'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0511'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('D51').Value = '0.0₇0972'
$ws.Range('E51').Value = '  -0.98%  '
